$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (OKLO)
$ws.Range("D2").Value = 108.78
$ws.Range("E2").Value = 55.5
$ws.Range("F2").Value = 22.61
$ws.Range("I2").Value = 73
$ws.Range("K2").Value = 59.9
$ws.Range("N2").Value = 52.47848103381103

# Row 3 (NuScale)
$ws.Range("D3").Value = 22.65
$ws.Range("E3").Value = 48.4
$ws.Range("F3").Value = 18.96
$ws.Range("N3").Value = 52.47848103381103
